$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking values
# (e.g. "1.000", "30.414.48") stay as text strings, matching the source
# workbook where these cells are inlineStr, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.414.48"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.869.07"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.83"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2916"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06491"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.20"
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.65"
$ws.Range("E12").Value = "  +3.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7368"
$ws.Range("E13").Value = "  +5.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.873.21"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.138"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "272.76"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.406.19"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.41"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007533"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.121.37"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.220"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.168"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.273"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.44"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.921"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1002"
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.503"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.275"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.124"
$ws.Range("E33").Value = "  +3.07%  "
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6929"
$ws.Range("E36").Value = "  +1.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9999"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.716"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("E39").Value = "  +0.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.741"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.305"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.71"
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.967"
$ws.Range("E43").Value = "  +4.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4187"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.8354"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.88"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.225"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.48"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.991"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "918.03"
$ws.Range("E51").Value = "  -1.73%  "
